$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0.2958904109589041
$ws.Cells.Item(2, 5).Value = 0.1972602739726027
$ws.Cells.Item(2, 6).Value = 0.2958904109589041
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 0.6904109589041096

$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0.2958904109589041
$ws.Cells.Item(3, 5).Value = 0.1972602739726027
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0.2958904109589041
$ws.Cells.Item(3, 8).Value = 0.1972602739726027
$ws.Cells.Item(3, 9).Value = 0.1972602739726027

$ws.Cells.Item(4, 2).Value = 0.2958904109589041
$ws.Cells.Item(4, 3).Value = 0.2958904109589041
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.09863013698630137
$ws.Cells.Item(4, 6).Value = 0.09863013698630137
$ws.Cells.Item(4, 7).Value = 0.1972602739726027
$ws.Cells.Item(4, 8).Value = 0.09863013698630137
$ws.Cells.Item(4, 9).Value = 0.2958904109589041

$ws.Cells.Item(5, 2).Value = 0.1972602739726027
$ws.Cells.Item(5, 3).Value = 0.1972602739726027
$ws.Cells.Item(5, 4).Value = 0.09863013698630137
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0.1972602739726027
$ws.Cells.Item(5, 7).Value = 0.1972602739726027
$ws.Cells.Item(5, 8).Value = 0.09863013698630137
$ws.Cells.Item(5, 9).Value = 0.2958904109589041

$ws.Cells.Item(6, 2).Value = 0.2958904109589041
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0.09863013698630137
$ws.Cells.Item(6, 5).Value = 0.1972602739726027
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0.2958904109589041
$ws.Cells.Item(6, 8).Value = 0.09863013698630137
$ws.Cells.Item(6, 9).Value = 0.3945205479452055

$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0.2958904109589041
$ws.Cells.Item(7, 4).Value = 0.1972602739726027
$ws.Cells.Item(7, 5).Value = 0.1972602739726027
$ws.Cells.Item(7, 6).Value = 0.2958904109589041
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0.6904109589041096

$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0.1972602739726027
$ws.Cells.Item(8, 4).Value = 0.09863013698630137
$ws.Cells.Item(8, 5).Value = 0.09863013698630137
$ws.Cells.Item(8, 6).Value = 0.09863013698630137
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0
$ws.Cells.Item(8, 9).Value = 0.09863013698630137

$ws.Cells.Item(9, 2).Value = 0.6904109589041096
$ws.Cells.Item(9, 3).Value = 0.1972602739726027
$ws.Cells.Item(9, 4).Value = 0.2958904109589041
$ws.Cells.Item(9, 5).Value = 0.2958904109589041
$ws.Cells.Item(9, 6).Value = 0.3945205479452055
$ws.Cells.Item(9, 7).Value = 0.6904109589041096
$ws.Cells.Item(9, 8).Value = 0.09863013698630137
$ws.Cells.Item(9, 9).Value = 0
